{"js": "// Implements \"Implemented Aug 18 2024 feedback, sans new sheriff info sheet\":\n//   1. Wrap \"Sign the petition.\" in a Jinja2 conditional\n//      (`{% if e_signature == False %} ... {% endif %}`) in the\n//      \"print/save the forms\" checklist item, and drop the now-redundant\n//      space before/after the two sentences that used to separate the runs.\n//   2. In the court-hearing checklist item, change \"court clerk\" to\n//      \"Circuit Clerk\" and fix the \"delivery\" -> \"deliver\" typo.\nconst body = context.document.body;\n\n// --- Edit 1: wrap \"Sign the petition.\" in an e_signature conditional -----\n// Original: \"...are correct. Sign the petition. You may need...\"\n// New:      \"...are correct.{% if e_signature == False %} Sign the\n//            petition.{% endif %} You may need...\"\nconst signStart = body.search(\"correct. Sign\", { matchCase: true, matchWholeWord: false });\nsignStart.load(\"items\");\nawait context.sync();\nif (signStart.items.length === 0) {\n  throw new Error(\"Could not find 'correct. Sign' text to update.\");\n}\nsignStart.items[0].insertText(\n  \"correct.{% if e_signature == False %} Sign\",\n  \"Replace\"\n);\nawait context.sync();\n\nconst signEnd = body.search(\"petition. You\", { matchCase: true, matchWholeWord: false });\nsignEnd.load(\"items\");\nawait context.sync();\nif (signEnd.items.length === 0) {\n  throw new Error(\"Could not find 'petition. You' text to update.\");\n}\nsignEnd.items[0].insertText(\n  \"petition.{% endif %} You\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Edit 2: \"court clerk\" -> \"Circuit Clerk\", \"delivery\" -> \"deliver\" ---\nconst clerkPhrase = body.search(\"Ask the court clerk\", { matchCase: true, matchWholeWord: false });\nclerkPhrase.load(\"items\");\nawait context.sync();\nif (clerkPhrase.items.length === 0) {\n  throw new Error(\"Could not find 'Ask the court clerk' text to update.\");\n}\nclerkPhrase.items[0].insertText(\"Ask the Circuit Clerk\", \"Replace\");\nawait context.sync();\n\nconst deliveryPhrase = body.search(\"need to delivery copies\", { matchCase: true, matchWholeWord: false });\ndeliveryPhrase.load(\"items\");\nawait context.sync();\nif (deliveryPhrase.items.length === 0) {\n  throw new Error(\"Could not find 'need to delivery copies' text to update.\");\n}\ndeliveryPhrase.items[0].insertText(\"need to deliver copies\", \"Replace\");\nawait context.sync();\n", "ps1": "# Implements \"Implemented Aug 18 2024 feedback, sans new sheriff info sheet\":\n#   1. Wrap \"Sign the petition.\" in a Jinja2 conditional\n#      (`{% if e_signature == False %} ... {% endif %}`) in the\n#      \"print/save the forms\" checklist item, and drop the now-redundant\n#      space before/after the two sentences that used to separate the runs.\n#   2. In the court-hearing checklist item, change \"court clerk\" to\n#      \"Circuit Clerk\" and fix the \"delivery\" -> \"deliver\" typo.\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($searchText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $searchText\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $found = $rng.Find.Execute()\n    if ($found) {\n        $rng.Text = $replaceText\n    } else {\n        Write-Output \"NOT FOUND: $searchText\"\n    }\n}\n\n# --- Edit 1: wrap \"Sign the petition.\" in an e_signature conditional -----\n# Original: \"...are correct. Sign the petition. You may need...\"\n# New:      \"...are correct.{% if e_signature == False %} Sign the\n#            petition.{% endif %} You may need...\"\nReplace-FirstMatch \"correct. Sign\" \"correct.{% if e_signature == False %} Sign\"\nReplace-FirstMatch \"petition. You\" \"petition.{% endif %} You\"\n\n# --- Edit 2: \"court clerk\" -> \"Circuit Clerk\", \"delivery\" -> \"deliver\" ---\nReplace-FirstMatch \"Ask the court clerk\" \"Ask the Circuit Clerk\"\nReplace-FirstMatch \"need to delivery copies\" \"need to deliver copies\"\n"}
